# Update loading_percent values for the 380 kV case.
# Rows 2-25 correspond to time steps 0-23; columns B, C, E, F, G, I, K, M, N
# hold the per-line loading percentages that were recomputed for this case
# (columns A, D, H, J, L, O are left untouched, matching the source diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("B","C","E","F","G","I","K","M","N")
$data = @(
    @(2, 11.52197107881983, 7.85231458517254, 11.02027059877949, 16.86991607391233, 3.689491207302202, 26.91105347852788, 12.2943842095911, 15.35652415228537, 21.27587702653777),
    @(3, 11.30779248317189, 7.678181610661778, 10.79380394866382, 15.89584955866808, 3.692260505092112, 26.94013856770037, 12.15505823811158, 15.21178828931428, 21.33209970627703),
    @(4, 11.17760832545187, 7.571872821007293, 10.65602210512897, 15.26997757108491, 3.694048494115865, 26.96311565800254, 12.07209165006378, 15.12648832516062, 21.36856137812239),
    @(5, 11.12496854336536, 7.528773226515678, 10.6002796073886, 15.00819731993403, 3.694799227551738, 26.97376282027511, 12.03896963683171, 15.09266059439035, 21.38390816670057),
    @(6, 11.1162548972665, 7.52163201761636, 10.5910505598934, 14.96433081551593, 3.694925224163837, 26.97560820811807, 12.03351238122484, 15.08710085297468, 21.38648599575648),
    @(7, 11.17689663726233, 7.571290573819334, 10.6552685902624, 15.26647399323137, 3.694058529129905, 26.96325405601053, 12.07164212108832, 15.12602829084732, 21.36876637249546),
    @(8, 11.44789499035941, 7.792187928728403, 10.94197244271261, 16.53996406344768, 3.690427921062203, 26.92001771838257, 12.24583279050389, 15.3059004158399, 21.29486002090195),
    @(9, 11.986386449984, 8.227254994693972, 11.51057006512644, 19.00274580682531, 3.684000025607744, 26.8759832126931, 12.60612599015936, 15.68529895415436, 21.16531208164947),
    @(10, 12.38183518823342, 8.544255207780807, 11.92738626003016, 20.67494806633232, 3.679694120027816, 26.86864381968296, 12.87980580682865, 15.97799363521163, 21.07948681079937),
    @(11, 12.56079055520717, 8.687146218035732, 12.11584167376793, 21.3917225636224, 3.677824656752574, 26.8707646318153, 13.00573153019776, 16.11366344369112, 21.04246907465074),
    @(12, 12.62835034331213, 8.741008788603425, 12.18696335874732, 21.65686569030329, 3.677129501265313, 26.87235433884311, 13.05357750309656, 16.16535663903434, 21.0287423235425),
    @(13, 12.61381043802779, 8.729420399545287, 12.17165797215958, 21.60004134736742, 3.677278648670486, 26.8719769630185, 13.04326657821545, 16.15421015450641, 21.03168568079756),
    @(14, 12.56635328293429, 8.69158279758091, 12.12169817368131, 21.4136618050453, 3.67776721038916, 26.87087964612352, 13.0096648818617, 16.11791015737984, 21.04133393358457),
    @(15, 12.53725540466683, 8.668372323300915, 12.09106261787886, 21.29868154950795, 3.678068129358369, 26.87030998438854, 12.98910247491254, 16.09571544409799, 21.0472816701593),
    @(16, 12.37011495103156, 8.53488544342769, 11.91504039334727, 20.62722412089977, 3.67981808478389, 26.86861521996246, 12.87160116411224, 15.96917422642568, 21.08194674025131),
    @(17, 12.26728749669095, 8.45261644158283, 11.80670488591016, 20.20408069597325, 3.68091444886555, 26.86897508043131, 12.79985164425851, 15.89215946708128, 21.10373112751402),
    @(18, 12.20805983996285, 8.405177062092042, 11.74428861564388, 19.95656407809801, 3.681553459591871, 26.86969587620887, 12.75872022732413, 15.84810264938532, 21.11645149927751),
    @(19, 12.18799419601829, 8.389096054220554, 11.72313996657343, 19.87204792380568, 3.681771264511632, 26.87002811968131, 12.74481873357121, 15.83322834879487, 21.12079113167426),
    @(20, 12.27824288445272, 8.461387014290015, 11.81824877033475, 20.24955283636154, 3.680796869053507, 26.86888358236856, 12.807475647701, 15.9003332993136, 21.10139241976369),
    @(21, 12.58029876702461, 8.702703771508862, 12.13637973786588, 21.46857628470577, 3.677623361989506, 26.87118059656005, 13.0195305084638, 16.12856407302871, 21.03849210896931),
    @(22, 12.77647393521067, 8.85895328432629, 12.34285226218178, 22.22866616901552, 3.67562368924788, 26.87726714076204, 13.15903658973087, 16.27955875047484, 20.99907963845333),
    @(23, 12.67190751263617, 8.775712278941667, 12.232810156967, 21.82633154458858, 3.67668416907983, 26.87359869062399, 13.08451037243298, 16.19881713697337, 21.01995960144587),
    @(24, 12.27329029314727, 8.457422275677482, 11.81303018789177, 20.22900810905287, 3.680849999790645, 26.86892334783402, 12.80402846368966, 15.89663722109496, 21.10244913816923),
    @(25, 11.84042531369147, 8.109757837521762, 11.35657476671213, 18.34778573295695, 3.68566540906445, 26.88351389960431, 12.50690311924498, 15.58003980001822, 21.19871402139793)
)
foreach ($row in $data) {
    $rowIndex = $row[0]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $cellRef = "$($columns[$i])$rowIndex"
        $ws.Range($cellRef).Value = $row[$i + 1]
    }
}
